# New crime data collected — weekly refresh of the 116th Precinct CompStat
# report: bumps the report "Volume/Number" and the covered week dates, and
# refreshes the crime-count table (Week to Date / 28 Day / Year to Date / 2
# Year columns C:L) for rows 15-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number + week-covered dates -----------------------------
$ws.Range("A8").Value  = "Volume 32   Number  17"
$ws.Range("C9").Value  = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# --- Row 15 (Rape) -----------------------------------------------------------
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("L15").Value = 100

# --- Row 16 (Robbery) --------------------------------------------------------
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("E16").Value = -100
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = -30.769230769230
$ws.Range("L16").Value = -33.333333333333

# --- Row 17 (Fel. Assault) ---------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -17.647058823529
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5.797101449275

# --- Row 18 (Burglary) -------------------------------------------------------
$ws.Range("C18").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").NumberFormat = $ws.Range("G18").NumberFormat
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -83.333333333333
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = -30
$ws.Range("L18").Value = -41.666666666666

# --- Row 19 (Gr. Larceny) ----------------------------------------------------
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 250
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 37.5
$ws.Range("I19").Value = 76
$ws.Range("J19").Value = 90
$ws.Range("K19").Value = -15.555555555555
$ws.Range("L19").Value = -13.636363636363

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 60
$ws.Range("L20").Value = 10.526315789473

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 75
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 12.727272727272
$ws.Range("I21").Value = 254
$ws.Range("J21").Value = 261
$ws.Range("K21").Value = -2.681992337164
$ws.Range("L21").Value = -5.925925925925

# --- Row 24 (Petit Larceny) --------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 44.444444444444
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = -19.148936170212
$ws.Range("I24").Value = 163
$ws.Range("J24").Value = 196
$ws.Range("K24").Value = -16.836734693877
$ws.Range("L24").Value = -8.938547486033

# --- Row 25 (Retail Theft) ---------------------------------------------------
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = -46.153846153846
$ws.Range("I25").Value = 37
$ws.Range("J25").Value = 38
$ws.Range("K25").Value = -2.631578947368
$ws.Range("L25").Value = 12.121212121212

# --- Row 26 (Misd. Assault) --------------------------------------------------
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 140
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 13.636363636363
$ws.Range("I26").Value = 116
$ws.Range("J26").Value = 126
$ws.Range("K26").Value = -7.936507936507
$ws.Range("L26").Value = 28.888888888888

# --- Row 27 (UCR Rape*) -------------------------------------------------------
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("L27").Value = 18.181818181818

# --- Row 28 (Other Sex Crimes) -----------------------------------------------
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -66.666666666666
